$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-only updates (column D)
$ws.Range("D2").Value = "'272.20"
$ws.Range("D3").Value = "'23.09"
$ws.Range("D4").Value = "'6.378"
$ws.Range("D5").Value = "'0.06276"
$ws.Range("D6").Value = "'3.644"
$ws.Range("D7").Value = "'6.738"
$ws.Range("D9").Value = "'0.8385"
$ws.Range("D10").Value = "'0.1632"
$ws.Range("D11").Value = "'0.08483"
$ws.Range("D12").Value = "'0.03491"
$ws.Range("D13").Value = "'0.03139"
$ws.Range("D40").Value = "'0.04695"
$ws.Range("D41").Value = "'0.006899"
$ws.Range("D42").Value = "'0.1176"
$ws.Range("D43").Value = "'0.003327"
$ws.Range("D44").Value = "'0.01254"
$ws.Range("D45").Value = "'0.00006251"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D47").Value = "'0.7975"
$ws.Range("D48").Value = "'0.09842"
$ws.Range("D49").Value = "'0.00001398"
$ws.Range("D50").Value = "'0.01238"

# Row shuffle rows 14-26: B (Coin), C (Link), D (Price), E (Volume label)
$ws.Range("B14").Value = "ProBitToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D14").Value = "'0.1263"
$ws.Range("E14").Value = "13ProBitTokenPROB"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.980"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").Value = "'0.09307"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001712"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04849"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006279"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.005482"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.001089"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.735"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.359"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D25").Value = "'0.01385"
$ws.Range("E25").Value = "24OneONE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3408"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
